$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H98").Value = 1496.5
$ws.Range("I98").Value = 1496.5
$ws.Range("K98").Value = 1496.5
$ws.Range("M98").Value = 1.5
$ws.Range("H99").Value = 209.57143
$ws.Range("I99").Value = 209.57143
$ws.Range("K99").Value = 628.71429
$ws.Range("M99").Value = 869.28571
$ws.Range("H122").Value = 1496.5
$ws.Range("I122").Value = 1496.5
$ws.Range("K122").Value = 4489.5
$ws.Range("M122").Value = -2039.5
$ws.Range("H125").Value = 7622.5
$ws.Range("I125").Value = 12245
$ws.Range("K125").Value = 110205
$ws.Range("M125").Value = -107745
$ws.Range("H132").Value = 1296.3636
$ws.Range("I132").Value = 940.5208
$ws.Range("K132").Value = 2821.5624
$ws.Range("M132").Value = -291.5623999999998

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 11217.069
$ws.Range("I32").Value = 9779.538
$ws.Range("K32").Value = 9779.538
$ws.Range("M32").Value = -9492.538
$ws.Range("H43").Value = 46093.715
$ws.Range("I43").Value = 41298.5
$ws.Range("J43").Value = 52487.332
$ws.Range("K43").Value = 41298.5
$ws.Range("L43").Value = 52487.332
$ws.Range("M43").Value = -40985.5
$ws.Range("N43").Value = -53113.332
$ws.Range("H45").Value = 122190.12
$ws.Range("I45").Value = 203321.9
$ws.Range("J45").Value = 6287.5713
$ws.Range("K45").Value = 203321.9
$ws.Range("L45").Value = 6287.5713
$ws.Range("M45").Value = -202944.9
$ws.Range("N45").Value = -7041.5713
$ws.Range("H61").Value = 8101.207
$ws.Range("I61").Value = 8862.458000000001
$ws.Range("K61").Value = 8862.458000000001
$ws.Range("M61").Value = -8650.458000000001
$ws.Range("H74").Value = 3198.3513
$ws.Range("I74").Value = 2588.516
$ws.Range("J74").Value = 6349.1665
$ws.Range("K74").Value = 2588.516
$ws.Range("L74").Value = 6349.1665
$ws.Range("M74").Value = -1714.516
$ws.Range("N74").Value = -8097.1665
$ws.Range("H77").Value = 3198.3513
$ws.Range("I77").Value = 2588.516
$ws.Range("J77").Value = 6349.1665
$ws.Range("K77").Value = 12942.58
$ws.Range("L77").Value = 31745.8325
$ws.Range("M77").Value = -8574.58
$ws.Range("N77").Value = -40481.8325
$ws.Range("H132").Value = 1819.2
$ws.Range("I132").Value = 1790.8334
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5372.5002
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2842.5002
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 8101.207
$ws.Range("I136").Value = 8862.458000000001
$ws.Range("K136").Value = 26587.374
$ws.Range("M136").Value = -24037.374

$ws = $wb.Worksheets("BSM")
$ws.Range("H82").Value = 20551.4
$ws.Range("J82").Value = 41750
$ws.Range("L82").Value = 41750
$ws.Range("N82").Value = -42516
$ws.Range("H85").Value = 20551.4
$ws.Range("J85").Value = 41750
$ws.Range("L85").Value = 41750
$ws.Range("N85").Value = -44402
$ws.Range("H94").Value = 2549.7778
$ws.Range("I94").Value = 1916
$ws.Range("J94").Value = 2866.6667
$ws.Range("K94").Value = 1916
$ws.Range("L94").Value = 2866.6667
$ws.Range("M94").Value = -1465
$ws.Range("N94").Value = -3768.6667
$ws.Range("H107").Value = 2424.25
$ws.Range("I107").Value = 2484.8572
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 2484.8572
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -564.8571999999999
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets("CRP")
$ws.Range("H51").Value = 10550
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4264
$ws.Range("H52").Value = 90389.164
$ws.Range("J52").Value = 90389.164
$ws.Range("L52").Value = 90389.164
$ws.Range("N52").Value = -90977.164
$ws.Range("H58").Value = 10717.45
$ws.Range("I58").Value = 4992.25
$ws.Range("K58").Value = 4992.25
$ws.Range("M58").Value = -4789.25
$ws.Range("H60").Value = 19030.334
$ws.Range("I60").Value = 15125
$ws.Range("K60").Value = 15125
$ws.Range("M60").Value = -14614
$ws.Range("H61").Value = 10550
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4652
$ws.Range("H69").Value = 21967.4
$ws.Range("I69").Value = 9999.5
$ws.Range("K69").Value = 9999.5
$ws.Range("M69").Value = -9250.5
$ws.Range("H72").Value = 21967.4
$ws.Range("I72").Value = 9999.5
$ws.Range("K72").Value = 29998.5
$ws.Range("M72").Value = -26254.5
$ws.Range("H107").Value = 1865.7142
$ws.Range("I107").Value = 567.75
$ws.Range("K107").Value = 567.75
$ws.Range("M107").Value = 1352.25
$ws.Range("H134").Value = 3808.1538
$ws.Range("I134").Value = 2174.6667
$ws.Range("K134").Value = 6524.000100000001
$ws.Range("M134").Value = -3989.000100000001
$ws.Range("H136").Value = 10717.45
$ws.Range("I136").Value = 4992.25
$ws.Range("K136").Value = 14976.75
$ws.Range("M136").Value = -12426.75

$ws = $wb.Worksheets("CUL")
$ws.Range("H70").Value = 1000
$ws.Range("I70").Value = 1000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2685
$ws.Range("H73").Value = 1000
$ws.Range("I73").Value = 1000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -1908
$ws.Range("H75").Value = 853.25
$ws.Range("I75").Value = 963
$ws.Range("K75").Value = 2889
$ws.Range("M75").Value = -1891
$ws.Range("H78").Value = 853.25
$ws.Range("I78").Value = 963
$ws.Range("K78").Value = 8667
$ws.Range("M78").Value = -3675

$ws = $wb.Worksheets("GSM")
$ws.Range("H99").Value = 5164.3
$ws.Range("I99").Value = 5164.3
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 5164.3
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2918.3
$ws.Range("N99").ClearContents()
$ws.Range("H114").Value = 54109.2
$ws.Range("J114").Value = 54109.2
$ws.Range("L114").Value = 54109.2
$ws.Range("N114").Value = -62787.2
$ws.Range("H122").Value = 1710
$ws.Range("I122").Value = 775
$ws.Range("K122").Value = 2325
$ws.Range("M122").Value = 125

$ws = $wb.Worksheets("LTW")
$ws.Range("H132").Value = 8889.394
$ws.Range("I132").Value = 9735.321
$ws.Range("K132").Value = 29205.963
$ws.Range("M132").Value = -26675.963

$ws = $wb.Worksheets("WVR")
$ws.Range("H100").Value = 1332
$ws.Range("I100").Value = 1531.4
$ws.Range("K100").Value = 3062.8
$ws.Range("M100").Value = -2521.8
$ws.Range("H122").Value = 2375.0232
$ws.Range("I122").Value = 1656.9706
$ws.Range("K122").Value = 4970.9118
$ws.Range("M122").Value = -2520.9118
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H132").Value = 3761.2026
$ws.Range("I132").Value = 3295.8032
$ws.Range("K132").Value = 9887.409599999999
$ws.Range("M132").Value = -7357.409599999999
$ws.Range("H136").Value = 3312
$ws.Range("I136").Value = 2685.9
$ws.Range("J136").Value = 6442.5
$ws.Range("K136").Value = 8057.700000000001
$ws.Range("L136").Value = 19327.5
$ws.Range("M136").Value = -5507.700000000001
$ws.Range("N136").Value = -24427.5
